$d = $word.ActiveDocument

function Insert-ParagraphAfterMatch($doc, $targetText, $newText) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs($i)
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $targetText) {
            $p.Range.InsertParagraphAfter()
            $doc.Paragraphs($i + 1).Range.Text = $newText
            return $true
        }
    }
    return $false
}

# Add "ID: L1" right after the "Name: WebSearch" paragraph
# (i.e. immediately before the "Description: WebSearch ..." paragraph).
Insert-ParagraphAfterMatch $d "Name: WebSearch" "ID: L1" | Out-Null

# Add "ID: L2" right after the "Name: GetHelp" paragraph
# (i.e. immediately before the "Description: Directs the user ..." paragraph).
Insert-ParagraphAfterMatch $d "Name: GetHelp" "ID: L2" | Out-Null
